$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each target cell originally stores its value as text (e.g. "299.71", "-4.82%"),
# not as a number, so force text formatting before assigning the new value --
# otherwise numeric-looking strings like "300.32" would be auto-converted to
# real numbers by Excel's normal cell-input parsing.
$updates = @(
    @{ Cell = "D2"; Value = "300.32" }
    @{ Cell = "E2"; Value = "-4.73%" }
    @{ Cell = "D3"; Value = "35.16" }
    @{ Cell = "E3"; Value = "-0.81%" }
    @{ Cell = "D4"; Value = "5.048" }
    @{ Cell = "E4"; Value = "-1.03%" }
    @{ Cell = "D5"; Value = "0.07934" }
    @{ Cell = "E5"; Value = "-2.88%" }
    @{ Cell = "D6"; Value = "1.907" }
    @{ Cell = "E6"; Value = "-7.47%" }
    @{ Cell = "D7"; Value = "7.781" }
    @{ Cell = "E7"; Value = "-2.08%" }
    @{ Cell = "B8"; Value = "GateToken" }
    @{ Cell = "C8"; Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt" }
    @{ Cell = "D8"; Value = "4.029" }
    @{ Cell = "E8"; Value = "-2.78%" }
    @{ Cell = "B9"; Value = "MXToken" }
    @{ Cell = "C9"; Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx" }
    @{ Cell = "D9"; Value = "0.9270" }
    @{ Cell = "E9"; Value = "-0.49%" }
    @{ Cell = "B10"; Value = "LiechtensteinCryptoassetsExchange" }
    @{ Cell = "C10"; Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx" }
    @{ Cell = "D10"; Value = "0.1328" }
    @{ Cell = "E10"; Value = "27.92%" }
    @{ Cell = "B11"; Value = "WazirX" }
    @{ Cell = "C11"; Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx" }
    @{ Cell = "D11"; Value = "0.1902" }
    @{ Cell = "E11"; Value = "-0.97%" }
    @{ Cell = "B12"; Value = "MandalaExchangeToken" }
    @{ Cell = "C12"; Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx" }
    @{ Cell = "D12"; Value = "0.09095" }
    @{ Cell = "E12"; Value = "0.25%" }
    @{ Cell = "B13"; Value = "BitrueCoin" }
    @{ Cell = "C13"; Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr" }
    @{ Cell = "D13"; Value = "0.03464" }
    @{ Cell = "E13"; Value = "-3.73%" }
    @{ Cell = "B14"; Value = "BitMartToken" }
    @{ Cell = "C14"; Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx" }
    @{ Cell = "D14"; Value = "0.09884" }
    @{ Cell = "E14"; Value = "-0.06%" }
    @{ Cell = "B15"; Value = "BitForexToken" }
    @{ Cell = "C15"; Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf" }
    @{ Cell = "D15"; Value = "0.001398" }
    @{ Cell = "E15"; Value = "-2.71%" }
    @{ Cell = "B16"; Value = "TigerCash" }
    @{ Cell = "C16"; Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch" }
    @{ Cell = "D16"; Value = "0.005720" }
    @{ Cell = "E16"; Value = "0.71%" }
    @{ Cell = "B17"; Value = "LEO" }
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo" }
    @{ Cell = "D17"; Value = "3.510" }
    @{ Cell = "E17"; Value = "1.21%" }
    @{ Cell = "D18"; Value = "2.962" }
    @{ Cell = "E18"; Value = "3.83%" }
    @{ Cell = "D19"; Value = "0.3406" }
    @{ Cell = "E19"; Value = "-0.03%" }
    @{ Cell = "D20"; Value = "0.1295" }
    @{ Cell = "E20"; Value = "-1.12%" }
    @{ Cell = "D21"; Value = "5.033" }
    @{ Cell = "E21"; Value = "-1.20%" }
    @{ Cell = "D22"; Value = "0.2404" }
    @{ Cell = "E22"; Value = "8.55%" }
    @{ Cell = "D23"; Value = "0.04497" }
    @{ Cell = "E23"; Value = "-1.08%" }
    @{ Cell = "D24"; Value = "0.001213" }
    @{ Cell = "E24"; Value = "-2.13%" }
    @{ Cell = "D25"; Value = "0.004757" }
    @{ Cell = "E25"; Value = "-0.81%" }
    @{ Cell = "D26"; Value = "0.0001230" }
    @{ Cell = "E26"; Value = "-1.59%" }
    @{ Cell = "D27"; Value = "0.0003000" }
    @{ Cell = "E27"; Value = "-33.33%" }
    @{ Cell = "D39"; Value = "0.01884" }
    @{ Cell = "E39"; Value = "-4.83%" }
    @{ Cell = "D40"; Value = "0.04703" }
    @{ Cell = "E40"; Value = "-5.47%" }
    @{ Cell = "D41"; Value = "0.007360" }
    @{ Cell = "E41"; Value = "-3.22%" }
    @{ Cell = "D42"; Value = "0.009909" }
    @{ Cell = "E42"; Value = "25.90%" }
    @{ Cell = "D43"; Value = "0.1317" }
    @{ Cell = "E43"; Value = "-4.71%" }
    @{ Cell = "D44"; Value = "0.002110" }
    @{ Cell = "E44"; Value = "-6.63%" }
    @{ Cell = "D45"; Value = "0.009340" }
    @{ Cell = "E45"; Value = "-20.64%" }
    @{ Cell = "D46"; Value = "0.00006256" }
    @{ Cell = "E46"; Value = "-5.56%" }
    @{ Cell = "D47"; Value = "0.00000000750" }
    @{ Cell = "E47"; Value = "0.00%" }
    @{ Cell = "D48"; Value = "64.75" }
    @{ Cell = "E48"; Value = "1.06%" }
    @{ Cell = "D49"; Value = "0.001659" }
    @{ Cell = "E49"; Value = "-2.43%" }
    @{ Cell = "D50"; Value = "0.00002100" }
    @{ Cell = "E50"; Value = "0.00%" }
    @{ Cell = "D51"; Value = "0.0002000" }
    @{ Cell = "E51"; Value = "0.00%" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}
